# Sync up for big compy
#
# 1. Remove the "CROSS OVER ON DP/DN" (red-highlighted, ilvl=1) bullet
#    from the U2/U3 placement list.
# 2. The document's stray "_GoBack" bookmark (left over from the author's
#    last cursor position) moves along with the edit: it now sits at the
#    start of the paragraph that used to follow the deleted bullet
#    ("Place C10, L2 and U4 near J1.") instead of its old spot on the
#    "Use figure below for example layout of U6." paragraph.

$d = $word.ActiveDocument

# --- Step 1: delete the "CROSS OVER ON DP/DN" paragraph entirely ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*CROSS OVER ON DP/DN*") {
        $p.Range.Delete()
        break
    }
}

# --- Step 2: relocate the "_GoBack" bookmark ---
# Remove it from its old home (the "Use figure below..." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create it, collapsed, at the start of the paragraph that now
# immediately follows the deleted bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Place C10, L2 and U4 near J1.*") {
        $target = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $target)
        break
    }
}
